# doubleClick.xlsx v0.4 restructure:
#   target=common options (columns A-D), value=options by case (columns E-G)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Propagate borders/alignment into the new columns (E:G) BEFORE touching
#    values, by copying formats from structurally-equivalent existing cells.
#    This re-uses existing style records instead of minting near-duplicates.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null          # xlPasteFormats

$ws.Range("D1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("F2:G2").PasteSpecial(-4122) | Out-Null

$ws.Range("A3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null

$ws.Range("D3").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Values - row 1 (headers): target = common options, value = per-case
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "open"
$ws.Range("C1").Value = "doubleClick"
$ws.Range("D1").Value = "sendKeys"
$ws.Range("E1").Value = "open"
$ws.Range("F1").Value = "doubleClick"
$ws.Range("G1").Value = "wait"

# ---------------------------------------------------------------------------
# 3. Values - row 2
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "http://127.0.0.1:9001/doubleclick/"
$ws.Range("C2").Value = "id=btn1"
$ws.Range("D2").Value = "xpath=//body"
$ws.Range("E2").Value = "http://127.0.0.1:9001/doubleclick/"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""

# ---------------------------------------------------------------------------
# 4. Values - row 3
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = '${ENTER_KEYS}'
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = '{"target": "id=btn1"}'
$ws.Range("G3").Value = 2000

# ---------------------------------------------------------------------------
# 5. Hyperlinks: re-point B2 at the new host, add a matching one on E2
# ---------------------------------------------------------------------------
$ws.Range("B2").Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "http://127.0.0.1:9001/doubleclick/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "http://127.0.0.1:9001/doubleclick/") | Out-Null

# ---------------------------------------------------------------------------
# 6. Font refresh: the v0.4 sheet swaps the old 12pt "Sarasa Mono CL" /
#    "ＭＳ Ｐゴシック" code-font runs for the workbook's own 11pt 游ゴシック,
#    keeping each run's original colour.
# ---------------------------------------------------------------------------
foreach ($addr in @("C1","D1","F1","G1","C2","F2","G2")) {
    $f = $ws.Range($addr).Font
    $f.Name = "游ゴシック"
    $f.Size = 11
    $f.Color = 526344   # FF080808
}

$f = $ws.Range("D2").Font
$f.Name = "游ゴシック"
$f.Size = 11
$f.Color = 1539863   # FF067D17 (green)

# ---------------------------------------------------------------------------
# 7. Rich text run colouring for F3: {"target": "id=btn1"}
#    { 1                default
#    "target" 2-9       purple FF871094
#    ": "     10-11     dark   FF080808
#    "id=btn1" 12-20    green  FF067D17
#    } 21               dark   FF080808
# ---------------------------------------------------------------------------
$cell = $ws.Range("F3")

$r = $cell.Characters(2, 8)
$r.Font.Name = "游ゴシック"; $r.Font.Size = 11; $r.Font.Color = 9703559

$r = $cell.Characters(10, 2)
$r.Font.Name = "游ゴシック"; $r.Font.Size = 11; $r.Font.Color = 526344

$r = $cell.Characters(12, 9)
$r.Font.Name = "游ゴシック"; $r.Font.Size = 11; $r.Font.Color = 1539334

$r = $cell.Characters(21, 1)
$r.Font.Name = "游ゴシック"; $r.Font.Size = 11; $r.Font.Color = 526344

# ---------------------------------------------------------------------------
# 8. Column widths (bestFit, character units) and final selection
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.7143
$ws.Columns.Item(2).ColumnWidth = 32.4286
$ws.Columns.Item(3).ColumnWidth = 10.2857
$ws.Columns.Item(4).ColumnWidth = 15.8571
$ws.Columns.Item(5).ColumnWidth = 32.4286
$ws.Columns.Item(6).ColumnWidth = 22.1429
$ws.Columns.Item(7).ColumnWidth = 4.7143

$ws.Range("E7").Select() | Out-Null

Write-Host "doubleClick.xlsx updated to v0.4 layout"
